$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking total (row 12 = "Total", B12 = corrects, E12 = "corr/total" summary)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
